$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H3 ("PlantUML status" column) from "Complete?" to "Complete" and apply the "Good" cell style
$ws.Range("H3").Value = "Complete"
$ws.Range("H3").Style = "Good"

# Update H9 ("PlantUML status" column) from "Changes done" to "Complete" and apply the "Good" cell style
$ws.Range("H9").Value = "Complete"
$ws.Range("H9").Style = "Good"

# Update the view: scroll so column E is the left-most visible column, and select I3
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("I3").Select()
